# Generate Report for Handback
#
# This applies the "handback" update to the localization-status workbook:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    for the two localizable files, on the Overview sheet and on each
#    per-language sheet.
#  - Each per-language sheet gains "Latest Target File" (E) and
#    "Latest Handback File" (F) hyperlinks (pointing at the same file the
#    row already references for handoff), and the
#    "Latest Handback DateTime" (G) is updated from the "never happened"
#    placeholder to the real handback timestamp.

$wb = $excel.ActiveWorkbook

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: update the Status cells for both languages.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $statusNew
$overview.Range("C2").Value = $statusNew
$overview.Range("B3").Value = $statusNew
$overview.Range("C3").Value = $statusNew

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("B2").Value = $statusNew
$zh.Range("B3").Value = $statusNew

$zhMdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/f18a230cdf4d8a7664962363a668879e57a62e49/e2e/84a055a9-c9f0-480a-a2f4-36cc40a6bfac.md"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/894c2b01249eab131ba1d974790eb7f787584a20/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/84a055a9-c9f0-480a-a2f4-36cc40a6bfac.0e695c2f58776d25a88612f2756d33575b4ca5af.zh-cn.xlf"
$mdDisplay = "84a055a9-c9f0-480a-a2f4-36cc40a6bfac.md"
$zhXlfDisplay = "84a055a9-c9f0-480a-a2f4-36cc40a6bfac.0e695c2f58776d25a88612f2756d33575b4ca5af.zh-cn.xlf"

$zh.Hyperlinks.Add($zh.Range("E2"), $zhMdTarget, "", "", $mdDisplay)
$zh.Hyperlinks.Add($zh.Range("F2"), $zhXlfTarget, "", "", $zhXlfDisplay)
$zh.Hyperlinks.Add($zh.Range("E3"), $zhMdTarget, "", "", $mdDisplay)
$zh.Hyperlinks.Add($zh.Range("F3"), $zhXlfTarget, "", "", $zhXlfDisplay)

$zh.Range("G2").Value = "2016-03-08 23:30:18"
$zh.Range("G3").Value = "2016-03-08 23:30:18"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("B2").Value = $statusNew
$de.Range("B3").Value = $statusNew

$deMdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/f18a230cdf4d8a7664962363a668879e57a62e49/e2e/84a055a9-c9f0-480a-a2f4-36cc40a6bfac.md"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0a414601ee6c7d264d815da96e152b3a03ca65f5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/84a055a9-c9f0-480a-a2f4-36cc40a6bfac.0e695c2f58776d25a88612f2756d33575b4ca5af.de-de.xlf"
$deXlfDisplay = "84a055a9-c9f0-480a-a2f4-36cc40a6bfac.0e695c2f58776d25a88612f2756d33575b4ca5af.de-de.xlf"

$de.Hyperlinks.Add($de.Range("E2"), $deMdTarget, "", "", $mdDisplay)
$de.Hyperlinks.Add($de.Range("F2"), $deXlfTarget, "", "", $deXlfDisplay)
$de.Hyperlinks.Add($de.Range("E3"), $deMdTarget, "", "", $mdDisplay)
$de.Hyperlinks.Add($de.Range("F3"), $deXlfTarget, "", "", $deXlfDisplay)

$de.Range("G2").Value = "2016-03-08 23:30:45"
$de.Range("G3").Value = "2016-03-08 23:30:45"

$wb.Save()
